# Auto-generated edit script: updates crypto price/volume table per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '243.57'
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.NumberFormat = "@"
$cell.Value = '-0.14%'
$cell.Style = "Normal"

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '29.81'
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.NumberFormat = "@"
$cell.Value = '13.71%'
$cell.Style = "Normal"

# Row 4
$cell = $ws.Range("E4")
$cell.NumberFormat = "@"
$cell.Value = '-0.27%'
$cell.Style = "Normal"

# Row 5
$cell = $ws.Range("E5")
$cell.NumberFormat = "@"
$cell.Value = '1.48%'
$cell.Style = "Normal"

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '6.525'
$cell.Style = "Normal"
$cell = $ws.Range("E6")
$cell.NumberFormat = "@"
$cell.Value = '0.81%'
$cell.Style = "Normal"

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.8411'
$cell.Style = "Normal"
$cell = $ws.Range("E7")
$cell.NumberFormat = "@"
$cell.Value = '2.35%'
$cell.Style = "Normal"

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.8649'
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.NumberFormat = "@"
$cell.Value = '3.32%'
$cell.Style = "Normal"

# Row 9
$cell = $ws.Range("B9")
$cell.NumberFormat = "@"
$cell.Value = 'WazirX'
$cell.Style = "Normal"
$cell = $ws.Range("C9")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$cell.Style = "Normal"
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.1340'
$cell.Style = "Normal"
$cell = $ws.Range("E9")
$cell.NumberFormat = "@"
$cell.Value = '0.89%'
$cell.Style = "Normal"

# Row 10
$cell = $ws.Range("B10")
$cell.NumberFormat = "@"
$cell.Value = 'MandalaExchangeToken'
$cell.Style = "Normal"
$cell = $ws.Range("C10")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$cell.Style = "Normal"
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.06915'
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.NumberFormat = "@"
$cell.Value = '-1.00%'
$cell.Style = "Normal"

# Row 11
$cell = $ws.Range("B11")
$cell.NumberFormat = "@"
$cell.Value = 'BitrueCoin'
$cell.Style = "Normal"
$cell = $ws.Range("C11")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$cell.Style = "Normal"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.02893'
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.NumberFormat = "@"
$cell.Value = '0.33%'
$cell.Style = "Normal"

# Row 12
$cell = $ws.Range("B12")
$cell.NumberFormat = "@"
$cell.Value = 'BitMartToken'
$cell.Style = "Normal"
$cell = $ws.Range("C12")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$cell.Style = "Normal"
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.09380'
$cell.Style = "Normal"
$cell = $ws.Range("E12")
$cell.NumberFormat = "@"
$cell.Value = '-0.13%'
$cell.Style = "Normal"

# Row 13
$cell = $ws.Range("B13")
$cell.NumberFormat = "@"
$cell.Value = 'BitForexToken'
$cell.Style = "Normal"
$cell = $ws.Range("C13")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '0.001521'
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.NumberFormat = "@"
$cell.Value = '0.29%'
$cell.Style = "Normal"

# Row 14
$cell = $ws.Range("B14")
$cell.NumberFormat = "@"
$cell.Value = 'CoinExToken'
$cell.Style = "Normal"
$cell = $ws.Range("C14")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '0.04162'
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.NumberFormat = "@"
$cell.Value = '-10.61%'
$cell.Style = "Normal"

# Row 15
$cell = $ws.Range("B15")
$cell.NumberFormat = "@"
$cell.Value = 'One'
$cell.Style = "Normal"
$cell = $ws.Range("C15")
$cell.NumberFormat = "@"
$cell.Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$cell.Style = "Normal"
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.0006019'
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.NumberFormat = "@"
$cell.Value = '-93.99%'
$cell.Style = "Normal"

# Row 16
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '0.005995'
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.NumberFormat = "@"
$cell.Value = '-4.13%'
$cell.Style = "Normal"

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '3.508'
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.NumberFormat = "@"
$cell.Value = '-3.83%'
$cell.Style = "Normal"

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.022'
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.NumberFormat = "@"
$cell.Value = '-0.37%'
$cell.Style = "Normal"

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '2.242'
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.NumberFormat = "@"
$cell.Value = '2.70%'
$cell.Style = "Normal"

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.3150'
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.NumberFormat = "@"
$cell.Value = '1.22%'
$cell.Style = "Normal"

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '0.03276'
$cell.Style = "Normal"
$cell = $ws.Range("E21")
$cell.NumberFormat = "@"
$cell.Value = '5.52%'
$cell.Style = "Normal"

# Row 22
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '0.1295'
$cell.Style = "Normal"
$cell = $ws.Range("E22")
$cell.NumberFormat = "@"
$cell.Value = '-0.30%'
$cell.Style = "Normal"

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '3.620'
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.NumberFormat = "@"
$cell.Value = '-3.14%'
$cell.Style = "Normal"

# Row 24
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '0.1374'
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.NumberFormat = "@"
$cell.Value = '-0.03%'
$cell.Style = "Normal"

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '0.001212'
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.NumberFormat = "@"
$cell.Value = '-2.84%'
$cell.Style = "Normal"

# Row 26
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '0.004312'
$cell.Style = "Normal"
$cell = $ws.Range("E26")
$cell.NumberFormat = "@"
$cell.Value = '-4.14%'
$cell.Style = "Normal"

# Row 27
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '0.0001180'
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.NumberFormat = "@"
$cell.Value = '22.90%'
$cell.Style = "Normal"

# Row 28
$cell = $ws.Range("E28")
$cell.NumberFormat = "@"
$cell.Value = '0.29%'
$cell.Style = "Normal"

# Row 40
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.03714'
$cell.Style = "Normal"
$cell = $ws.Range("E40")
$cell.NumberFormat = "@"
$cell.Value = '2.10%'
$cell.Style = "Normal"

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.005330'
$cell.Style = "Normal"
$cell = $ws.Range("E41")
$cell.NumberFormat = "@"
$cell.Value = '-13.20%'
$cell.Style = "Normal"

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.1058'
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.NumberFormat = "@"
$cell.Value = '0.65%'
$cell.Style = "Normal"

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '0.002311'
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.NumberFormat = "@"
$cell.Value = '-3.72%'
$cell.Style = "Normal"

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '0.009799'
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.NumberFormat = "@"
$cell.Value = '10.49%'
$cell.Style = "Normal"

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.00005112'
$cell.Style = "Normal"
$cell = $ws.Range("E45")
$cell.NumberFormat = "@"
$cell.Value = '-4.61%'
$cell.Style = "Normal"

# Row 46
$cell = $ws.Range("E46")
$cell.NumberFormat = "@"
$cell.Value = '-0.01%'
$cell.Style = "Normal"

# Row 47
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.09997'
$cell.Style = "Normal"
$cell = $ws.Range("E47")
$cell.NumberFormat = "@"
$cell.Value = '-30.58%'
$cell.Style = "Normal"

# Row 48
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.002725'
$cell.Style = "Normal"
$cell = $ws.Range("E48")
$cell.NumberFormat = "@"
$cell.Value = '18.71%'
$cell.Style = "Normal"

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '0.00002100'
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.NumberFormat = "@"
$cell.Value = '-0.01%'
$cell.Style = "Normal"

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '0.0002000'
$cell.Style = "Normal"
$cell = $ws.Range("E50")
$cell.NumberFormat = "@"
$cell.Value = '-0.01%'
$cell.Style = "Normal"
